$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '22.020.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -1.83%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'" + '1.555.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -1.00%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'" + '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.13%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'" + '1.0000'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -0.08%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'" + '286.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.37%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = "'" + '0.3782'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +2.50%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'" + '0.3232'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -2.47%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('B9').Value = "'" + 'OKB'
$ws.Range('B9').Style = 'Normal'
$ws.Range('C9').Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C9').Style = 'Normal'
$ws.Range('D9').Value = "'" + '41.19'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -13.22%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('B10').Value = "'" + 'Polygon'
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').Value = "'" + 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').Value = "'" + '1.124'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -3.79%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'" + '0.07299'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -2.64%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = "'" + '1.0000'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -0.12%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'" + '19.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -6.61%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'" + '5.717'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -3.57%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'" + '6.771'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -1.74%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = "'" + '1.553.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -0.32%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'" + '0.00001085'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -2.54%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'" + '0.06628'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -1.55%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'" + '85.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -3.45%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'" + '6.422'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -0.12%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = "'" + '0.9990'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.12%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'" + '15.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -3.25%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'" + '11.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -4.34%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'" + '22.033.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -1.70%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = "'" + '2.275'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -4.05%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'" + '2.527'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -3.59%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = "'" + '148.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -1.70%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'" + '18.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -3.80%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'" + '4.859'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -1.41%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = "'" + '1.727.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -0.64%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'" + '119.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -3.84%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'" + '1.114'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +3.06%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'" + '5.918'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -2.71%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'" + '1.653'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -16.38%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('B35').Value = "'" + 'FraxShare'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'" + 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'" + '9.268'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -6.10%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('B36').Value = "'" + 'Stellar'
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.08144'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -2.07%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('B37').Value = "'" + 'InternetComputer(DFINITY)'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'" + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'" + '5.235'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -1.57%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').Value = "'" + 'Hedera'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'" + '0.06179'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -3.11%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = "'" + '0.02279'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -6.29%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'" + '0.2113'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -4.22%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Value = "'" + '1.216'
$ws.Range('D41').Style = 'Normal'

$ws.Range('D42').Value = "'" + '10.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -4.32%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = "'" + '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -0.01%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'" + '0.5935'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -4.61%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'" + '13.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -3.08%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'" + '3.717'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -1.54%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'" + '0.5725'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -5.07%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').Value = "'" + 'NEARProtocol'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'" + '1.931'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -5.08%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = "'" + 'Quant'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'" + '119.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -4.03%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = "'" + '1.156'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -3.16%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = "'" + '0.06888'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -4.22%  '
$ws.Range('E51').Style = 'Normal'

